$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-13 01:51:29"

for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
